$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trial")
$ws.Activate()

$ws.Range("C1").Value = 75
$ws.Range("C2").Value = -1147.5

$ws.Range("C8").Select()
